$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "progresivas" table grew by one more milestone row (row 70). Fill in
# the new record: a kilometric marker ("Hito 01") with no width/length
# measurements ("-").
$ws.Range("A70").Value = 11070
$ws.Range("B70").Value = "-"
$ws.Range("C70").Value = "-"
$ws.Range("D70").Value = "Hito_kilometrico"
$ws.Range("E70").Value = "Hito 01"
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 1

# Column F ("progresivas_al_reves") is the distance back from the last
# milestone. Re-point every existing row at the new last row (A70) and
# give the freshly appended row its own copy of the same formula.
$ws.Range("F2:F69").Formula = "=`$A`$70-A2"
$ws.Range("F70").Formula = "=`$A`$70-A70"

# Mirror the author's final selection/cursor position.
$ws.Range("I70").Select()
